$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 0.0005555555555555556
$ws.Range("K2").Value = 3909
$ws.Range("L2").Value = 0.007818
